$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codes")

# --- Update existing cells (Code column corrections, typo fixes, new Survey Code column) ---
$ws.Range("D1").Value = "Survey Code"
$ws.Range("C2").Value = "XFT-SM1"
$ws.Range("D2").Value = "X1"
$ws.Range("C3").Value = "XFT-1"
$ws.Range("D3").Value = "X2"
$ws.Range("C4").Value = "XFT-2"
$ws.Range("D4").Value = "X3"
$ws.Range("C5").Value = "XFT-3"
$ws.Range("D5").Value = "X4"
$ws.Range("C6").Value = "XFT-4"
$ws.Range("D6").Value = "X5"
$ws.Range("C7").Value = "XFT-5"
$ws.Range("D7").Value = "X6"
$ws.Range("C8").Value = "XFT-PG1"
$ws.Range("D8").Value = "X7"
$ws.Range("C9").Value = "OPO1"
$ws.Range("D9").Value = "X8"
$ws.Range("C10").Value = "DM1"
$ws.Range("D10").Value = "X9"
$ws.Range("C11").Value = "PgM1"
$ws.Range("D11").Value = "X10"
$ws.Range("C13").Value = "Dsgnr1"
$ws.Range("B15").Value = "CCB-adm."
$ws.Range("C15").Value = "CCB-A"
$ws.Range("C16").Value = "TC1"
$ws.Range("C19").Value = "CAKM1"
$ws.Range("A21").Value = "Helena Eberil"
$ws.Range("C22").Value = "Dsgnr2"

# --- New rows appended to the Codes table ---
$ws.Range("A27").Value = "Lars Rundberg"
$ws.Range("B27").Value = "Design"
$ws.Range("C27").Value = "Dsgnr3"
$ws.Range("A28").Value = "Eva Cullman"
$ws.Range("C28").Value = "O7"
$ws.Range("A29").Value = "Mikael Krekola"
$ws.Range("B29").Value = "PG-froCpp"
$ws.Range("C29").Value = "PG3"
$ws.Range("A30").Value = "Pierre Svärd"
$ws.Range("B30").Value = "Designer"
$ws.Range("C30").Value = "Dsgnr4"
$ws.Range("A31").Value = "Anny Lei"
$ws.Range("B31").Value = "Feature Proj. Leader"
$ws.Range("C31").Value = "FPjL1"
$ws.Range("A32").Value = "Jan Johansson"
$ws.Range("B32").Value = "Prod. Manager RBS Sys"
$ws.Range("C32").Value = "PdMRBS1"
$ws.Range("A33").Value = "Per Simonsson"
$ws.Range("B33").Value = "RBS System"
$ws.Range("C33").Value = "RBS1"
$ws.Range("A34").Value = "Henrik Sundh"
$ws.Range("B34").Value = "Strategic Product Manager"
$ws.Range("C34").Value = "SPM1"
$ws.Range("A35").Value = "Jeanette Munro"
$ws.Range("B35").Value = "Designer"
$ws.Range("C35").Value = "Dsgnr5"
$ws.Range("A36").Value = "Ricardo Morales"
$ws.Range("B36").Value = "Designer"
$ws.Range("C36").Value = "Dsgnr6"
$ws.Range("A37").Value = "Carl Ohlsson"
$ws.Range("B37").Value = "Designer"
$ws.Range("C37").Value = "Dsgnr7"
$ws.Range("A38").Value = "Niklas Isaksson"
$ws.Range("B38").Value = "Program Manager"
$ws.Range("C38").Value = "PgM3"
$ws.Range("A39").Value = "Thomas Andersson"
$ws.Range("B39").Value = "Sector Manager"
$ws.Range("C39").Value = "SrM1"
$ws.Range("A40").Value = "Per Lofter"
$ws.Range("B40").Value = "Section Manager"
$ws.Range("C40").Value = "SM3"
$ws.Range("A41").Value = "Thomas Nyberg"
$ws.Range("B41").Value = "OPO/LC team"
$ws.Range("C41").Value = "OPO4"

# --- Column widths (bestFit-style autosize for the new/extended columns) ---
$ws.Columns.Item(2).ColumnWidth = 21.9166666666667
$ws.Columns.Item(4).ColumnWidth = 10.6666666666667

# --- View state: Codes sheet becomes the active/selected tab, Day-1 loses it ---
$ws.Activate()
$ws.Range("C41").Select()
